# Updates cryptos list (price/volume refresh), per GitHub Actions run.
# Note: some "Price" values are single-decimal numeric-looking strings
# (e.g. "1.002"); a leading apostrophe is used to force Excel to store
# them as text (matching the original cell's text type) instead of
# auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.600.49'
$ws.Range("E2").Value = '  +1.40%  '
$ws.Range("D3").Value = '1.831.55'
$ws.Range("E3").Value = '  +1.43%  '
$ws.Range("D4").Value = '''1.002'
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = '''317.94'
$ws.Range("E5").Value = '  +0.17%  '
$ws.Range("E6").Value = '  +0.15%  '
$ws.Range("E7").Value = '  +0.66%  '
$ws.Range("D8").Value = '''0.3980'
$ws.Range("E8").Value = '  +5.20%  '
$ws.Range("D9").Value = '''0.07717'
$ws.Range("E9").Value = '  +2.99%  '
$ws.Range("D10").Value = '''1.120'
$ws.Range("D11").Value = '''41.91'
$ws.Range("E11").Value = '  -0.21%  '
$ws.Range("D12").Value = '''6.348'
$ws.Range("E12").Value = '  +2.01%  '
$ws.Range("D13").Value = '''20.98'
$ws.Range("E13").Value = '  +1.88%  '
$ws.Range("D14").Value = '''7.598'
$ws.Range("E14").Value = '  +3.20%  '
$ws.Range("E15").Value = '  +0.24%  '
$ws.Range("D16").Value = '1.828.89'
$ws.Range("E16").Value = '  +1.44%  '
$ws.Range("D17").Value = '''92.84'
$ws.Range("E17").Value = '  +3.47%  '
$ws.Range("E18").Value = '  +1.65%  '
$ws.Range("D19").Value = '''0.06570'
$ws.Range("E19").Value = '  +0.97%  '
$ws.Range("D20").Value = '''17.79'
$ws.Range("E20").Value = '  +3.06%  '
$ws.Range("D21").Value = '''1.002'
$ws.Range("E21").Value = '  +0.18%  '
$ws.Range("D22").Value = '''6.084'
$ws.Range("E22").Value = '  +2.68%  '
$ws.Range("D23").Value = '28.612.98'
$ws.Range("E23").Value = '  +1.36%  '
$ws.Range("D24").Value = '''11.21'
$ws.Range("D25").Value = '''2.243'
$ws.Range("E25").Value = '  +7.30%  '
$ws.Range("E26").Value = '  +1.21%  '
$ws.Range("D27").Value = '2.041.90'
$ws.Range("E27").Value = '  +1.50%  '
$ws.Range("D28").Value = '''156.17'
$ws.Range("E28").Value = '  +0.39%  '
$ws.Range("E29").Value = '  +3.91%  '
$ws.Range("D30").Value = '''125.10'
$ws.Range("E30").Value = '  +2.44%  '
$ws.Range("D31").Value = '''1.137'
$ws.Range("E31").Value = '  +1.08%  '
$ws.Range("D32").Value = '''0.1118'
$ws.Range("D33").Value = '''5.749'
$ws.Range("E33").Value = '  +2.79%  '
$ws.Range("D34").Value = '''3.659'
$ws.Range("E34").Value = '  +0.88%  '
$ws.Range("D35").Value = '''0.07248'
$ws.Range("E35").Value = '  +0.20%  '
$ws.Range("D36").Value = '''0.2256'
$ws.Range("E36").Value = '  +1.27%  '
$ws.Range("D37").Value = '''0.02351'
$ws.Range("E37").Value = '  +2.28%  '
$ws.Range("D38").Value = '''8.905'
$ws.Range("E38").Value = '  +4.55%  '
$ws.Range("D39").Value = '''5.210'
$ws.Range("E39").Value = '  +2.15%  '
$ws.Range("D40").Value = '''11.40'
$ws.Range("E40").Value = '  +2.37%  '
$ws.Range("D41").Value = '''0.6308'
$ws.Range("E41").Value = '  +2.14%  '
$ws.Range("D42").Value = '''1.198'
$ws.Range("E42").Value = '  +1.15%  '
$ws.Range("E43").Value = '  +0.21%  '
$ws.Range("D44").Value = '''1.395'
$ws.Range("E44").Value = '  -2.67%  '
$ws.Range("D45").Value = '''13.49'
$ws.Range("E45").Value = '  +1.19%  '
$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D46").Value = '''3.721'
$ws.Range("E46").Value = '  +1.07%  '
$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").Value = '''0.5901'
$ws.Range("E47").Value = '  +2.19%  '
$ws.Range("D48").Value = '''125.19'
$ws.Range("E48").Value = '  -0.56%  '
$ws.Range("E49").Value = '  +3.57%  '
$ws.Range("D50").Value = '''1.196'
$ws.Range("E50").Value = '  +0.08%  '
$ws.Range("D51").Value = '''0.06941'
$ws.Range("E51").Value = '  +1.78%  '
